# Applies the "Add files via upload" commit:
#   - Sheet renamed from IClientBalance-20241128-085807- to
#     IClientBalance-20241129-090503- (refreshed export timestamp)
#   - Every "Dt. Referencia" date in column G (rows 2-274) bumped by one
#     day: 45624 (2024-11-28) -> 45625 (2024-11-29)
#   - Row 224's "Saldo Previsto" (E) and "Vl. Total" (H) values updated
#     from 560.11 to 591.01

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the refreshed export timestamp.
$ws.Name = "IClientBalance-20241129-090503-"

# Bump every reference date in column G (rows 2 through 274) by one day.
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45625
}

# Row 224 ("E224"/"H224") balance correction.
$ws.Range("E224").Value = 591.01
$ws.Range("H224").Value = 591.01
